$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update header row labels (B2: Funcionalidade -> Requisito, D2: Prioridade -> Classificação)
# Set D2 first, then B2, so new shared-string entries are appended in the same
# order as the target workbook ("Classificação" before "Requisito").
$ws.Range("D2").Value = "Classificação"
$ws.Range("B2").Value = "Requisito"

# Update the view: zoom out to 60% and move the selection to L6
$excel.ActiveWindow.Zoom = 60
$ws.Range("L6").Select()
